$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently holds 5 data rows:
#   2: contact_email      | contact@datannur.com
#   3: alias_1             | institution : owner
#   4: alias_2             | institution : manager
#   5: filter_1             | open_data : Open Data
#   6: filter_2             | closed_data : Closed Data
#
# Internalize the institution alias rows (alias_1 / alias_2) by deleting
# them outright; the filter_1 / filter_2 rows shift up to become rows 3/4.
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()

# contact_email's value cell no longer carries the extra fill style.
$ws.Range("B2").Style = "Normal"

# Match the saved selection state from the edit.
$ws.Range("A6").Select()
